$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.100.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.857.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.84%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.60%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4377"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.64%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3670"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07463"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9305"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.19"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.868.03"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.659"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.398"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06910"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.64%  "
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "81.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008972"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.134.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.086"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.117.76"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.007"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.35%  "
$ws.Range("E27").Value = "  -3.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.272"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.716"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08974"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.813"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7879"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.168"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.974"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.004"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.122"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05402"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01958"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.948"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5235"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.969"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1668"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.652"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.06716"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4838"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "106.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.005"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.37%  "
$ws.Range("E50").Value = "  -7.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.663"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.75%  "
